$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scripts")

# D1: MaxEpochs 30 -> 31
$ws.Range("D1").Value = 31

# Extend the shared formulas in column E and F down to row 60
$ws.Range("E60").Formula = '=$D$3&" Both "&$D$1&" "&$D$2&"Client.xml"&" "&$D$2&"DataShape"&$C60&".xml "&$D$2&"trainDataSet"&$D60&".xml "&$D$2&"Engine"&$B60&".xml"'
$ws.Range("F60").Formula = '=$D$3&" Infer "&$D$1&" "&$D$2&"Client.xml"&" "&$D$2&"trainDataSet"&$D60&".xml "&$A60'

# New data row 60: Engine4, DataShape1, trainDataSet7
$ws.Range("B60").Value = 4
$ws.Range("C60").Value = 1
$ws.Range("D60").Value = 7

# Update the view: scroll the frozen pane & selection to match new extent
$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Range("E60").Select()

$wb.Save()
